$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 1000.6667
$ws.Range("I82").Value = 1000.6667
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 3002.0001
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -2596.0001
$ws.Range("N82").Value = ""

$ws.Range("H85").Value = 1000.6667
$ws.Range("I85").Value = 1000.6667
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 3002.0001
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -1598.0001
$ws.Range("N85").Value = ""

$ws.Range("H132").Value = 409169.9
$ws.Range("I132").Value = 434646.1
$ws.Range("J132").Value = 52503
$ws.Range("K132").Value = 1303938.3
$ws.Range("L132").Value = 157509
$ws.Range("M132").Value = -1301408.3
$ws.Range("N132").Value = -162569

$ws.Range("H137").Value = 90911320
$ws.Range("I137").Value = 250002430
$ws.Range("J137").Value = 2114
$ws.Range("K137").Value = 750007290
$ws.Range("L137").Value = 6342
$ws.Range("M137").Value = -750004740
$ws.Range("N137").Value = -11442

$ws.Range("H140").Value = 61612.5
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 61612.5
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 61612.5
$ws.Range("N140").Value = -71972.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 33750
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 33750
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 33750
$ws.Range("N7").Value = -33978

$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").Value = ""

$ws.Range("H74").Value = 5303.241
$ws.Range("I74").Value = 855.86365
$ws.Range("J74").Value = 19280.715
$ws.Range("K74").Value = 855.86365
$ws.Range("L74").Value = 19280.715
$ws.Range("M74").Value = 18.13634999999999
$ws.Range("N74").Value = -21028.715

$ws.Range("H77").Value = 5303.241
$ws.Range("I77").Value = 855.86365
$ws.Range("J77").Value = 19280.715
$ws.Range("K77").Value = 4279.31825
$ws.Range("L77").Value = 96403.575
$ws.Range("M77").Value = 88.68174999999974
$ws.Range("N77").Value = -105139.575

$ws.Range("H110").Value = 616.8889
$ws.Range("I110").Value = 519
$ws.Range("J110").Value = 1400
$ws.Range("K110").Value = 519
$ws.Range("L110").Value = 1400
$ws.Range("M110").Value = 1526

$ws.Range("H122").Value = 2513.5676
$ws.Range("I122").Value = 1775.871
$ws.Range("J122").Value = 6325
$ws.Range("K122").Value = 5327.613
$ws.Range("L122").Value = 18975
$ws.Range("M122").Value = -2877.613
$ws.Range("N122").Value = -23875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 381.66666
$ws.Range("I22").Value = 350
$ws.Range("J22").Value = 388
$ws.Range("K22").Value = 350
$ws.Range("L22").Value = 388
$ws.Range("M22").Value = -177

$ws.Range("H134").Value = 2166.7556
$ws.Range("I134").Value = 1432.1082
$ws.Range("J134").Value = 5564.5
$ws.Range("K134").Value = 4296.3246
$ws.Range("L134").Value = 16693.5
$ws.Range("M134").Value = -1761.3246
$ws.Range("N134").Value = -21763.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5109.1943
$ws.Range("I31").Value = 1856.238
$ws.Range("J31").Value = 9663.333000000001
$ws.Range("K31").Value = 1856.238
$ws.Range("L31").Value = 9663.333000000001
$ws.Range("M31").Value = -1561.238
$ws.Range("N31").Value = -10253.333

$ws.Range("H34").Value = 5109.1943
$ws.Range("I34").Value = 1856.238
$ws.Range("J34").Value = 9663.333000000001
$ws.Range("K34").Value = 1856.238
$ws.Range("L34").Value = 9663.333000000001
$ws.Range("M34").Value = -1654.238
$ws.Range("N34").Value = -10067.333

$ws.Range("H107").Value = 361.23077
$ws.Range("I107").Value = 210.11111
$ws.Range("J107").Value = 701.25
$ws.Range("K107").Value = 210.11111
$ws.Range("L107").Value = 701.25
$ws.Range("M107").Value = 1709.88889
$ws.Range("N107").Value = -4541.25

$ws.Range("H109").Value = 33599.8
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 33599.8
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 33599.8
$ws.Range("N109").Value = -35679.8

$ws.Range("H122").Value = 1182.8125
$ws.Range("I122").Value = 893.1539
$ws.Range("J122").Value = 2438
$ws.Range("K122").Value = 2679.4617
$ws.Range("L122").Value = 7314
$ws.Range("M122").Value = -229.4616999999998
$ws.Range("N122").Value = -12214

$ws.Range("H132").Value = 2243.102
$ws.Range("I132").Value = 2063.4358
$ws.Range("J132").Value = 2943.8
$ws.Range("K132").Value = 6190.307400000001
$ws.Range("L132").Value = 8831.400000000001
$ws.Range("M132").Value = -3660.307400000001
$ws.Range("N132").Value = -13891.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1072.909
$ws.Range("I92").Value = 1401
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 4203
$ws.Range("L92").Value = 3000
$ws.Range("M92").Value = -2955
$ws.Range("N92").Value = -5496

$ws.Range("H113").Value = 16667349
$ws.Range("I113").Value = 704.3125
$ws.Range("J113").Value = 35714944
$ws.Range("K113").Value = 2112.9375
$ws.Range("L113").Value = 107144832
$ws.Range("M113").Value = 57.0625
$ws.Range("N113").Value = -107149172

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 30000
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 30000
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 30000
$ws.Range("N42").Value = -30970

$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = ""
$ws.Range("N107").Value = ""

$ws.Range("H113").Value = 3477.5
$ws.Range("I113").Value = 1955.5
$ws.Range("J113").Value = 4999.5
$ws.Range("K113").Value = 1955.5
$ws.Range("L113").Value = 4999.5
$ws.Range("M113").Value = 214.5
$ws.Range("N113").Value = -9339.5

$ws.Range("H115").Value = 30000
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 30000
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 30000
$ws.Range("N115").Value = -32350

$ws.Range("H122").Value = 696852.9399999999
$ws.Range("I122").Value = 1112611.9
$ws.Range("J122").Value = 3921.3333
$ws.Range("K122").Value = 3337835.7
$ws.Range("L122").Value = 11763.9999
$ws.Range("M122").Value = -3335385.7
$ws.Range("N122").Value = -16663.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 6251009.5
$ws.Range("I16").Value = 7692934.5
$ws.Range("J16").Value = 2666.3333
$ws.Range("K16").Value = 7692934.5
$ws.Range("L16").Value = 2666.3333
$ws.Range("M16").Value = -7692764.5
$ws.Range("N16").Value = -3006.3333

$ws.Range("H93").Value = 1454
$ws.Range("I93").Value = 1230.6
$ws.Range("J93").Value = 1677.4
$ws.Range("K93").Value = 1230.6
$ws.Range("L93").Value = 1677.4
$ws.Range("M93").Value = 17.40000000000009
$ws.Range("N93").Value = -4173.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 33500
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 33500
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 33500
$ws.Range("N27").Value = -33638

$ws.Range("H115").Value = 26994
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 26994
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 26994
$ws.Range("N115").Value = -30128

$ws.Range("H132").Value = 3572.0286
$ws.Range("I132").Value = 3726.2593
$ws.Range("J132").Value = 3051.5
$ws.Range("K132").Value = 11178.7779
$ws.Range("L132").Value = 9154.5
$ws.Range("M132").Value = -8648.777900000001
$ws.Range("N132").Value = -14214.5

$ws.Range("H136").Value = 3423.5881
$ws.Range("I136").Value = 1840.1333
$ws.Range("J136").Value = 4673.684
$ws.Range("K136").Value = 5520.3999
$ws.Range("L136").Value = 14021.052
$ws.Range("M136").Value = -2970.3999
$ws.Range("N136").Value = -19121.052
